{"js": "// The document has a \"_GoBack\" bookmark (Word's \"last edit location\" marker)\n// sitting in the second paragraph, right before the final \":\" run. Word moved\n// it to the very start of the document (beginning of the first paragraph)\n// when the file was reopened/resaved. Reproduce that by deleting the\n// existing bookmark and re-inserting it, collapsed, at the start of the body.\n\ncontext.document.deleteBookmark(\"_GoBack\");\n\nconst start = context.document.body.getRange(\"Start\");\nstart.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# The \"_GoBack\" bookmark (Word's \"last edit location\" marker) currently\n# sits in the second paragraph, right before the final \":\" run. Reopening /\n# resaving the document in Word moved it to the very beginning of the\n# document (start of the first paragraph). Reproduce that: drop the old\n# bookmark and add a new, collapsed one at the start of the document.\n\n$d = $word.ActiveDocument\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Anchor a collapsed bookmark exactly at document position 0. A bookmark\n# range that is already collapsed at position 0 when handed to\n# Bookmarks.Add ends up split across the paragraph boundary, so instead we\n# bracket a transient character at the very start, bookmark that (non-\n# collapsed) range, and then delete the transient character again - the\n# bookmark collapses back down to position 0 and stays put.\n$head = $d.Range(0, 0)\n$head.InsertBefore(\"X\")\n\n$span = $d.Range(0, 1)\n$d.Bookmarks.Add(\"_GoBack\", $span)\n\n$d.Range(0, 1).Delete()\n"}
